$wb = $excel.ActiveWorkbook

# NOTE: the worksheet that carries the data (tab shown as selected, with the
# TSLA calcs) is the one named "Sheet2" in the workbook's sheet collection
# (the sheet named "Sheet1" is the other, essentially-empty sheet).
$ws = $wb.Worksheets.Item("Sheet2")

# Update the ML "Label" input column (C) - these are the cells whose raw
# inputs were adjusted; the dependent D/E/H/I/J formula cells recalculate
# automatically from these.
$ws.Range("C4").Value2 = 0
$ws.Range("C5").Value2 = 1
$ws.Range("C6").Value2 = 0
$ws.Range("C7").Value2 = 0
$ws.Range("C11").Value2 = 0
$ws.Range("C15").Value2 = 1
$ws.Range("C16").Value2 = 0
$ws.Range("C18").Value2 = 1
$ws.Range("C20").Value2 = 1
$ws.Range("C21").Value2 = 0
$ws.Range("C24").Value2 = 0
$ws.Range("C25").Value2 = 1
$ws.Range("C26").Value2 = 0
$ws.Range("C28").Value2 = 1
$ws.Range("C29").Value2 = 1
$ws.Range("C31").Value2 = 0
$ws.Range("C33").Value2 = 1
$ws.Range("C34").Value2 = 0
$ws.Range("C35").Value2 = 1
$ws.Range("C36").Value2 = 0
$ws.Range("C40").Value2 = 1
$ws.Range("C43").Value2 = 0
$ws.Range("C44").Value2 = 0
$ws.Range("C45").Value2 = 0
$ws.Range("C47").Value2 = 1
$ws.Range("C49").Value2 = 1
$ws.Range("C51").Value2 = 1

# Move the active selection, matching the author's final cursor position.
$ws.Activate()
$ws.Range("K6").Select()

# Best-effort: reflect the author's window position at save time.
$win = $excel.ActiveWindow
$win.Left = 5220
$win.Top = 4791
